# Update the threshold table values on Sheet1 and move the active
# selection, matching the authored edit (Zn/3His/Threshold -> Cu/3His/
# probability re-export): the alpha_distance_range (row 2) and
# beta_distance_range (row 3) Min/Max values were refreshed, and the
# last-saved selection moved to C5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (alpha_distance_range): Min 5.8 -> 5.5, Max 10.7 -> 11
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 11

# Row 3 (beta_distance_range): Min 5.7 -> 5.5, Max 9.3 -> 9.5
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 9.5

# Active cell/selection moved from C3 to C5 before the file was saved.
$ws.Range("C5").Select()
